# 4.0.3 model and data
# Expand the single consolidated "trans/BVTQaZ/BVTQaZ.csv" and
# "trans/VTQaZ/VTQaZ.csv" rows on the "Boolean" sheet into six per-mode
# CSV rows each (LDVs, HDVs, aircraft, rail, ships, motorbikes), and
# update the UI selection / active-tab state to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Boolean" sheet: split the two consolidated rows into six rows each
# ---------------------------------------------------------------------
$wsBool = $wb.Worksheets.Item("Boolean")

# --- trans/BVTQaZ/BVTQaZ.csv (row 17): replace with 6 per-mode rows.
#     Insert 5 extra rows below row 17 to make room.
$wsBool.Range("A18:A22").EntireRow.Insert()

$wsBool.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBool.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBool.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBool.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBool.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBool.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# --- trans/VTQaZ/VTQaZ.csv (originally row 21, now shifted to row 26
#     after the BVTQaZ insert above): replace with 6 per-mode rows.
$wsBool.Range("A27:A31").EntireRow.Insert()

$wsBool.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBool.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBool.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBool.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBool.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBool.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# --- six trailing blank (but explicitly formatted) rows at the bottom
#     of the used range, rows 33-38 -- touch their formatting (matching
#     the rest of column A) so they materialise as styled empty rows and
#     extend the sheet's used range/dimension down to row 38.
$wsBool.Range("A33:A38").Font.Name = "Calibri"
$wsBool.Range("A33:A38").Font.Size = 11

# Scroll/selection state for this sheet
$wsBool.Activate()
$wsBool.Range("A32").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. "Integer" sheet: selection moves to A13, tab selection moves away
# ---------------------------------------------------------------------
$wsInt = $wb.Worksheets.Item("Integer")
$wsInt.Activate()
$wsInt.Range("A13").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. "About" sheet becomes the active/selected tab
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
